# "explore options on format" - extend the GL export sheet with a mirrored
# block (module/gl/reason columns + a debit/credit pair per line item) and a
# running total, dropping the old stray "stuff"/"account" labels.
#
# Cell writes below are deliberately sequenced (module -> gl -> MJE_O ->
# reason -> supplies -> total -> the VBA note) to match the shared-string
# insertion order the authoring tool produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$itemNames = @{5="Green Chili"; 6="Black Beans"; 7="Distilled Water"; 8="Fruit Preservative"; 9="Watch Battery"; 10="Sales Tax"}
$amounts   = @{5=1.49; 6=1.6; 7=7.12; 8=3.99; 9=3.79; 10=0.26}
$reasons   = @{5="food"; 6="food"; 7="food"; 8="food"; 9="supplies"; 10="taxes"}

# 1) module
$ws.Range("D1").Value = "module"

# 2) gl
$ws.Range("E4").Value = "gl"

# 3) MJE_O
$ws.Range("D2").Value = "MJE_O"

# 4) reason (also used at G4)
$ws.Range("C4").Value = "reason"
$ws.Range("G4").Value = "reason"

# give C4/G4 (and the now-blank D4) the date-format style already used by
# B4/F4 ("s=1"), without minting a new style entry
$ws.Range("F4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").ClearContents()
$excel.CutCopyMode = 0

$ws.Range("H4").Value = "item"

# 5) supplies (first hit on row 9) - walk rows 5-10 filling C/D/E/F/G/H
foreach ($r in 5..10) {
    $reason = $reasons[$r]
    $ws.Range("C$r").Value = $reason
    $ws.Range("D$r").ClearContents()
    $ws.Range("E$r").Value = $reason
    $ws.Range("F$r").Value = $amounts[$r]
    $ws.Range("G$r").Value = $reason
    $ws.Range("H$r").Value = $itemNames[$r]
}

# 6) total
$ws.Range("F1").Value = "total"
$ws.Range("E1").Value = "offset"
$ws.Range("D2").Value = "MJE_O"
$ws.Range("E2").Value = "dcard"
$ws.Range("F2").Formula = "=SUM(B5:B10)"

# new offsetting (credit) lines, one per item row above
$offsetRows = @{11=5; 12=6; 13=7; 14=8; 15=9; 16=10}
foreach ($r in 11..16) {
    $src = $offsetRows[$r]
    $ws.Range("E$r").Value = "dcard"
    $ws.Range("F$r").Formula = "=-F$src"
    $ws.Range("G$r").Value = $reasons[$src]
    $ws.Range("H$r").Value = $itemNames[$src]
}

# 7) the explanatory note
$ws.Range("A13").Value = "VBA function will create full JSON with GL array"

# ---------------------------------------------------------------------
# Column widths for the newly-populated B/C columns, and move the
# active selection the way the author left it (D4).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

$ws.Range("D4").Select()
